$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "bleu" label with "noir" in column B (statut_label)
$ws.Columns.Item(2).Replace("bleu", "noir", 1, 1, $false, $false, $false, $false)

# Replace the status_name descriptions in column C (statut_name)
$ws.Columns.Item(3).Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", 1, 1, $false, $false, $false, $false)
$ws.Columns.Item(3).Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", 1, 1, $false, $false, $false, $false)
$ws.Columns.Item(3).Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", 1, 1, $false, $false, $false, $false)
$ws.Columns.Item(3).Replace("résultat et / ou publication posté", "résultat postés ou publiés", 1, 1, $false, $false, $false, $false)
